$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename (row 1)
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'
$ws.Range('B21').Value = 'Amatenango De La Frontera'
$ws.Range('B24').Value = 'Comitán De Domínguez'
$ws.Range('B37').Value = 'San Cristóbal De Las Casas'
$ws.Range('B55').Value = 'Hidalgo Del Parral'
$ws.Range('B76').Value = 'Villa De Álvarez'
$ws.Range('A78').Value = 'Ciudad De México'
$ws.Range('B82').Value = 'Cuajimalpa De Morelos'
$ws.Range('B97').Value = 'Nombre De Dios'
$ws.Range('B102').Value = 'San Juan De Guadalupe'
$ws.Range('B103').Value = 'San Juan Del Río'
$ws.Range('A108').Value = 'Estado De México'
$ws.Range('B108').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B110').Value = 'Almoloya De Juárez'
$ws.Range('B119').Value = 'Ecatepec De Morelos'
$ws.Range('B122').Value = 'Ixtapan De La Sal'
$ws.Range('B127').Value = 'Naucalpan De Juárez'
$ws.Range('B131').Value = 'San Felipe Del Progreso'
$ws.Range('B132').Value = 'Soyaniquilpan De Juárez'
$ws.Range('B136').Value = 'Tenango Del Valle'
$ws.Range('B143').Value = 'Tlalnepantla De Baz'
$ws.Range('B147').Value = 'Villa De Allende'
$ws.Range('A151').Value = 'Guanajuato'
$ws.Range('B154').Value = 'Apaseo El Alto'
$ws.Range('B155').Value = 'Apaseo El Grande'
$ws.Range('B159').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B163').Value = 'Jaral Del Progreso'
$ws.Range('B173').Value = 'San Francisco Del Rincón'
$ws.Range('B174').Value = 'Silao De La Victoria'
$ws.Range('B176').Value = 'Valle De Santiago'
$ws.Range('B179').Value = 'Acapulco De Juárez'
$ws.Range('B185').Value = 'Ayutla De Los Libres'
$ws.Range('B188').Value = 'Chilapa De Álvarez'
$ws.Range('B189').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B190').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B193').Value = 'Coyuca De Catalán'
$ws.Range('B195').Value = 'Iguala De La Independencia'
$ws.Range('B196').Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range('B197').Value = 'Zihuatanejo De Azueta'
$ws.Range('B198').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B206').Value = 'Taxco De Alarcón'
$ws.Range('B207').Value = 'Técpan De Galeana'
$ws.Range('B209').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B211').Value = 'Tixtla De Guerero'
$ws.Range('B213').Value = 'Tlalixtaquilla De Maldonado'
$ws.Range('B218').Value = 'Atotonilco El Grande'
$ws.Range('B222').Value = 'Huejutla De Reyes'
$ws.Range('B225').Value = 'Pachuca De Soto'
$ws.Range('B228').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B229').Value = 'Tezontepec De Aldama'
$ws.Range('B231').Value = 'Tula De Allende'
$ws.Range('B232').Value = 'Tulancingo De Bravo'
$ws.Range('B233').Value = 'Zacualtipán De Ángeles'
$ws.Range('B236').Value = 'Acatlán De Juárez'
$ws.Range('B237').Value = 'Ahualulco De Mercado'
$ws.Range('B241').Value = 'Atemajac De Brizuela'
$ws.Range('B242').Value = 'Atotonilco El Alto'
$ws.Range('B244').Value = 'Autlán De Navarro'
$ws.Range('B255').Value = 'Huejuquilla El Alto'
$ws.Range('B261').Value = 'Lagos De Moreno'
$ws.Range('B264').Value = 'Ojuelos De Jalisco'
$ws.Range('B267').Value = 'San Juanito De Escobedo'
$ws.Range('B269').Value = 'San Martín De Bolaños'
$ws.Range('B271').Value = 'San Miguel El Alto'
$ws.Range('B272').Value = 'San Sebastián Del Oeste'
$ws.Range('B274').Value = 'Tamazula De Gordiano'
$ws.Range('B276').Value = 'Teocuitatlán De Corona'
$ws.Range('B282').Value = 'Yahualica De González Gallo'
$ws.Range('B283').Value = 'Zacoalco De Torres'
$ws.Range('B285').Value = 'Zapotitlán De Vadillo'
$ws.Range('B286').Value = 'Zapotlán Del Rey'
$ws.Range('B287').Value = 'Zapotlán El Grande'
$ws.Range('B305').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B360').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B381').Value = 'Puente De Ixtla'
$ws.Range('B383').Value = 'Tlaltizapán De Zapata'
$ws.Range('B392').Value = 'Amatlán De Cañas'
$ws.Range('B393').Value = 'Bahía De Banderas'
$ws.Range('B396').Value = 'Ixtlán Del Río'
$ws.Range('B403').Value = 'Santa María Del Oro'
$ws.Range('B412').Value = 'Montemorelos'
$ws.Range('B414').Value = 'San Nicolás De Los Garza'
$ws.Range('B418').Value = 'Chalcatongo De Hidalgo'
$ws.Range('B420').Value = 'Coicoyán De Las Flores'
$ws.Range('B421').Value = 'Constancia Del Rosario'
$ws.Range('B423').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B424').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B425').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B426').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B427').Value = 'Mariscala De Juárez'
$ws.Range('B428').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B429').Value = 'Oaxaca De Juárez'
$ws.Range('B430').Value = 'Ocotlán De Morelos'
$ws.Range('B431').Value = 'Putla Villa De Guerero'
$ws.Range('B455').Value = 'Santo Domingo De Morelos'
$ws.Range('B458').Value = 'Teotitlán De Flores Magón'
$ws.Range('B460').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B461').Value = 'Zapotitlán Del Río'
$ws.Range('B462').Value = 'Zimatlán De Álvarez'
$ws.Range('B470').Value = 'Cuayuca De Andrade'
$ws.Range('B476').Value = 'Izúcar De Matamoros'
$ws.Range('B478').Value = 'Los Reyes De Juárez'
$ws.Range('B480').Value = 'Palmar De Bravo'
$ws.Range('B497').Value = 'Cadereyta De Montes'
$ws.Range('B500').Value = 'Jalpan De Serra'
$ws.Range('B502').Value = 'Pinal De Amoles'
$ws.Range('B505').Value = 'San Juan Del Río'
$ws.Range('B508').Value = 'Armadillo De Los Infante'
$ws.Range('B509').Value = 'Ciudad Del Maíz'
$ws.Range('B512').Value = 'Mexquitic De Carmona'
$ws.Range('B516').Value = 'Villa De Ramos'
$ws.Range('B565').Value = 'Muñoz De Domingo Arenas'
$ws.Range('B567').Value = 'Tetla De La Solidaridad'
$ws.Range('B572').Value = 'Amatlán De Los Reyes'
$ws.Range('B574').Value = 'Castillo De Teayo'
$ws.Range('B578').Value = 'Cosamaloapan De Carpio'
$ws.Range('B583').Value = 'Lerdo De Tejada'
$ws.Range('B585').Value = 'Martínez De La Torre'
$ws.Range('B588').Value = 'Paso De Ovejas'
$ws.Range('B589').Value = 'Poza Rica De Hidalgo'
$ws.Range('B592').Value = 'Sayula De Alemán'
$ws.Range('B610').Value = 'Concepción Del Oro'
$ws.Range('B620').Value = 'Mezquital Del Oro'
$ws.Range('B622').Value = 'Nochistlán De Mejía'
$ws.Range('B623').Value = 'Noria De Ángeles'
$ws.Range('B631').Value = 'Teúl De González Ortega'
$ws.Range('B632').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B634').Value = 'Villa De Cos'

# Numeric precision fix
$ws.Range('D327').Value = 0.009954058192955587

# Remove trailing footer rows 641-646 (also updates dimension to A1:D640)
$ws.Range("A641:D646").EntireRow.Delete()